$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D22: add the missing checkmark (Wingdings "ü"), matching the style used
# by every other filled-in cell in column D (e.g. D3, D20, D21, D23).
$ws.Range("D22").Value = "ü"
$ws.Range("D22").Font.Name = "Wingdings"
$ws.Range("D22").Font.Size = 12

# C24: add the missing checkmark (Wingdings "ü"), matching the style used
# by every other filled-in cell in column C (e.g. C3, C20, C21, C22).
$ws.Range("C24").Value = "ü"
$ws.Range("C24").Font.Name = "Wingdings"
$ws.Range("C24").Font.Size = 12
